# "att dicionario de dados"
# Fill in the previously-blank "Volume esperado:", "Tempo de retenção:" and
# "Rotina de limpeza:" rows on the "DD entidade" sheet with their actual
# values for each entity, and update the view/selection state to match.

$wb = $excel.ActiveWorkbook
$wsEntidade = $wb.Worksheets.Item("DD entidade")
$wsAtributo = $wb.Worksheets.Item("DD atributo")

# --- entidade empresa (rows 6-8) ---
$wsEntidade.Range("A6").Value  = "Volume esperado: 20 por mês"
$wsEntidade.Range("A7").Value  = "Tempo de retenção: permanente"
$wsEntidade.Range("A8").Value  = "Rotina de limpeza: Não se aplica"

# --- entidade caixa (rows 15-17) ---
$wsEntidade.Range("A15").Value = "Volume esperado: 200 por mês"
$wsEntidade.Range("A16").Value = "Tempo de retenção: 1 mês"
$wsEntidade.Range("A17").Value = "Rotina de limpeza:  Uma vez a cada 5 anos"

# --- entidade leitura (rows 24-26) ---
$wsEntidade.Range("A24").Value = "Volume esperado: 2000 por mês"
$wsEntidade.Range("A25").Value = "Tempo de retenção: 1 mês"
$wsEntidade.Range("A26").Value = "Rotina de limpeza: Uma vez a cada 5 anos"

# --- entidade orgao (rows 33-35) ---
$wsEntidade.Range("A33").Value = "Volume esperado: 20 por ano"
$wsEntidade.Range("A34").Value = "Tempo de retenção: permanente"
$wsEntidade.Range("A35").Value = "Rotina de limpeza: Uma vez a cada 1 ano"

# --- entidade rota (rows 42-44) ---
$wsEntidade.Range("A42").Value = "Volume esperado: 400 por mês"
$wsEntidade.Range("A43").Value = "Tempo de retenção: 1 mês"
$wsEntidade.Range("A44").Value = "Rotina de limpeza: Uma vez a cada 5 anos"

# --- entidade sensor (rows 51-53) ---
$wsEntidade.Range("A51").Value = "Volume esperado: 2000 por mês"
$wsEntidade.Range("A52").Value = "Tempo de retenção: 1 mês"
$wsEntidade.Range("A53").Value = "Rotina de limpeza: Uma vez a cada 5 anos"

# --- entidade usuario (row 60 only; 61/62 remain blank labels) ---
$wsEntidade.Range("A60").Value = "Volume esperado: carga inicial de 60 ocorrências e volume mensal de 2 ocorrências"

# --- restore view/selection state ---
# "DD atributo" was left scrolled back to the top with A56 selected.
$wsAtributo.Activate()
$wsAtributo.Range("A56").Select()

# "DD entidade" is the active tab, scrolled down with A60 selected.
$wsEntidade.Activate()
$wsEntidade.Range("A60").Select()
